$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correcting data analysis: use the fixed day-0 average (0.1904) as the
# normalisation denominator instead of chaining off the previous day's
# computed value, for C14 and for the shared formula spanning C15:C22.
$ws.Range("C14").Formula = "=(B14/0.1904*C13)"
$ws.Range("C15:C22").Formula = "=(B15/0.1904*C14)"

# Daily entry: update the active selection to the single cell C22.
$ws.Range("C22").Select()
